$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# Insert a new header/group row above row 3, shifting everything down.
$ws.Rows("3:3").Insert(-4121)   # xlShiftDown

$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"

# Give the new row similar formatting to the row below it (the real header row, now row 4).
$ws.Range("B3:E3").Borders.Item(9).LineStyle = 1   # xlEdgeBottom = 9, xlContinuous = 1
$ws.Range("B3:D3").Font.Name = $ws.Range("B4").Font.Name
$ws.Range("B3:D3").Font.Size = $ws.Range("B4").Font.Size

# Append a new closing row at the bottom of the table (row 20).
$ws.Range("B20:E20").Borders.Item(10).LineStyle = 1  # xlEdgeTop = 10

$wb.Save()
